$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Liens Ressources Humaines / K'IAM" rows (A4:D4, A16:D16,
# A27:D27, A37:D37) - clearing their contents drops the now-unreferenced
# shared strings ("Liens Ressources Humaines :", "K'IAM",
# "https://iamuaprod.kiabi.fr/idmdash", "Risorse Umane :") from the workbook.
$ws.Range("A4:D4").ClearContents()
$ws.Range("A16:D16").ClearContents()
$ws.Range("A27:D27").ClearContents()
$ws.Range("A37:D37").ClearContents()

# Add a hyperlink on D5 (3xONEY row) pointing at the Oney login page, same
# target used elsewhere in the sheet for the same link text.
$ws.Hyperlinks.Add($ws.Range("D5"), "https://open.oney.fr/", "/login")

# Update the active selection to match the saved state of the workbook.
[void]$ws.Range("D37").Select()
